$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.090.10"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.518.91"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "3.518.31"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -2.14%  "
$ws.Range("D13").Value = "4.117.51"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000179"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.25%  "
$ws.Range("D17").Value = "3.522.05"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "64.146.16"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  -2.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("E21").Value = "  -0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.92%  "
$ws.Range("D23").Value = "3.660.71"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").Value = "3.533.16"
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.57"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("E37").Value = "  -0.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "159.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.88%  "
$ws.Range("E42").Value = "  -2.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.812"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("E47").Value = "  -2.97%  "
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("D50").Value = "2.482.17"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("E51").Value = "  -0.86%  "
